$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.983.80'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '1.831.33'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9969'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.86'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6264'
$ws.Range("E6").Value = '  -4.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9977'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07586'
$ws.Range("E8").Value = '  +3.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2919'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07713'
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.836.53'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.947'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6645'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001026'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.83'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.049'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '28.995.69'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.03'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.33'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9968'
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.181'
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9977'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.20'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.491'
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1371'
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.87'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.489'
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.011'
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.187'
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05224'
$ws.Range("E32").Value = '  -2.41%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7349'
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.690'
$ws.Range("E36").Value = '  +1.74%  '
$ws.Range("D37").Value = '1.236.91'
$ws.Range("E37").Value = '  -4.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.751'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01784'
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.356'
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8951'
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9980'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.64'
$ws.Range("D44").Value = '1.982.02'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("E45").Value = '  +1.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.96'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5096'
$ws.Range("E47").Value = '  -0.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4038'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.865'
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05751'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.687'
$ws.Range("E51").Value = '  -0.47%  '
